# "update cummins program commit." — add the missing App-date entry (F4)
# for the row already describing the Stadium job, re-point the active
# selection at it, and give it the same yyyy/m date display used
# elsewhere on the sheet (no leading-quote text prefix this time, since
# it's a real date serial and not typed-as-text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data cell: F4 = 2012-10-12 (serial 41194).
$ws.Range("F4").Value = 41194
$ws.Range("F4").NumberFormat = 'yyyy"年"m"月"'

# The active cell in the frozen bottom-right pane moves from G5 to F5.
$ws.Range("F5").Select()
